# Natmi following Dr Hou advice
# Sheet1 previously held one L-R (C1qa-Cd93) result per Sending/Target
# cluster pair with "M2" as the only sending cluster. The update adds a
# second sending cluster ("sCs") and recomputes the stats for every
# Sending x Target cluster combination, growing the table from 4 to 8
# data rows (plus the header row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columnLetters = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# Each entry: row, Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# then numeric columns E..T in order.
$rowData = @(
  ,(2,"M2","C1qa","Cd93","ECs",3,1,191.5877736666667,574.763321,0.9786143588951871,0.9786143588951871,3,1,169.629438,508.888314,0.7428377317484701,0.7428377317484702,32498.92637474787,292490.3373727308,0.7269516706181841,0.7269516706181842)
  ,(3,"M2","C1qa","Cd93","FAPs",3,1,191.5877736666667,574.763321,0.9786143588951871,0.9786143588951871,3,1,0.9848756666666668,2.954627,0.004312947180081616,0.004312947180081616,188.6901363151408,1698.211226836267,0.004220712039584375,0.004220712039584375)
  ,(4,"M2","C1qa","Cd93","M2",3,1,191.5877736666667,574.763321,0.9786143588951871,0.9786143588951871,3,1,54.620752,163.862256,0.2391940691454494,0.2391940691454494,10464.66827167913,94182.01444511217,0.234078750628305,0.234078750628305)
  ,(5,"M2","C1qa","Cd93","sCs",3,1,191.5877736666667,574.763321,0.9786143588951871,0.9786143588951871,3,1,3.118221666666667,9.354665000000001,0.01365525192599884,0.01365525192599884,597.4131469158294,5376.718322242466,0.01336322560911363,0.01336322560911363)
  ,(6,"sCs","C1qa","Cd93","ECs",1,0.3333333333333333,4.186764,12.560292,0.02138564110481286,0.02138564110481286,3,1,169.629438,508.888314,0.7428377317484701,0.7428377317484702,710.198424358632,6391.785819227688,0.01588606113028603,0.01588606113028604)
  ,(7,"sCs","C1qa","Cd93","FAPs",1,0.3333333333333333,4.186764,12.560292,0.02138564110481286,0.02138564110481286,3,1,0.9848756666666668,2.954627,0.004312947180081616,0.004312947180081616,4.123441985676001,37.110977871084,0.00009223514049724012,0.00009223514049724012)
  ,(8,"sCs","C1qa","Cd93","M2",1,0.3333333333333333,4.186764,12.560292,0.02138564110481286,0.02138564110481286,3,1,54.620752,163.862256,0.2391940691454494,0.2391940691454494,228.684198126528,2058.157783138752,0.005115318517144372,0.005115318517144372)
  ,(9,"sCs","C1qa","Cd93","sCs",1,0.3333333333333333,4.186764,12.560292,0.02138564110481286,0.02138564110481286,3,1,3.118221666666667,9.354665000000001,0.01365525192599884,0.01365525192599884,13.05525821802,117.49732396218,0.0002920263168852158,0.0002920263168852158)
)

# Clear the previous 4-row table body first (A2:T5) so stale cells beyond
# the new table don't linger, then write the new 8-row body (A2:T9).
$ws.Range("A2:T5").Clear()

foreach ($entry in $rowData) {
    $r = $entry[0]
    for ($i = 1; $i -lt $entry.Length; $i++) {
        $colLetter = $columnLetters[$i - 1]
        $ws.Range("$colLetter$r").Value = $entry[$i]
    }
}
